$d = $word.ActiveDocument

# Move to the end of the document and insert a new paragraph
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Move(4, 1)  # wdParagraph, move forward by 1 paragraph

$end.Text = "Kwdsaihnjfpioehnfgpñioehnbgfpiewbndfpvinbedpif oeifjmnoeiphnfpñieowdnfbh"
